$d = $word.ActiveDocument

# 1) Update the delivery date text.
$d.Content.Find.Execute("Primera entrega 25/6/2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Primera entrega 26/6/2019", 2) | Out-Null

# 2) Move the hidden "_GoBack" bookmark from inside the "Proy01003/" run group
#    up to the very start of the "Ruta en GitLab:" paragraph (collapsed/empty bookmark
#    right after the paragraph's pPr, before its first run). Adding a bookmark with the
#    same name re-seats the existing one, so the old bookmarkStart/bookmarkEnd pair
#    disappears automatically.
$rutaPara = $d.Paragraphs.Item(3)
$startRng = $d.Range($rutaPara.Range.Start, $rutaPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $startRng) | Out-Null

# 3) Merge the split "Proy0100" + "3" + "/" runs (previously separated by the bookmark)
#    into a single run reading "Proy01003/". A direct Range.Text assignment to the exact
#    same text is treated as a no-op, so first swap in a scratch placeholder, then
#    replace that placeholder with the real text.
$target = $d.Content
$target.Find.Execute("Proy01003/") | Out-Null
$mergeRng = $d.Range($target.Start, $target.End)
$mergeRng.Text = "ZZPLACEHOLDERZZ"

$d.Content.Find.Execute("ZZPLACEHOLDERZZ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Proy01003/", 2) | Out-Null
